# Edit script: rewrite the Greenland Ice Sheet essay into a Governance/
# Democracy essay, per the commit diff.

$d = $word.ActiveDocument

function ReplaceText($old, $new) {
    $ok = $d.Content.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Host "WARNING: replace failed for: $old"
    }
}

# --- Title ---
ReplaceText "Unraveling the Secrets of the Greenland Ice Sheet" "Unraveling the Tapestry of Governance: The Interwoven Threads of Democracy"

# --- Author byline: "Dr. Elsie Robinson" -> "Emily Blackwood" (collapses to one run) ---
ReplaceText "Dr. Elsie Robinson" "Emily Blackwood"

# --- Email local part ---
ReplaceText "elsirobinson@glaciology" "emilyblackwood@educatormail"

# --- New body content blocks (vertical-tab char = <w:br/> line break) ---
$vtab = [char]11

$block1New = 'In the vast tapestry of human existence, there exists a realm where power and influence interlace to shape the destinies of nations and communities: the realm of governance. Governance, like a intricate web, comprises intricate threads of democracy, intertwining and shaping the fabric of our societies. Within this realm, citizens are bestowed with the sacred duty of participating in the decision-making processes that affect their lives. Our journey begins with an exploration of democracy, the bedrock upon which modern governance systems are built. Democracy, an intricate mosaic of individual voices, is a testament to the inherent power of collective decision-making. Through democratic processes, citizens have the opportunity to shape policies, elect representatives, and hold those in power accountable. However, democracy is not without its challenges. The delicate balance between majority rule and the protection of minority rights requires constant vigilance and unwavering commitment to principles of justice and equality. Time and time again, history has borne witness to the fragility of democracy, as populist movements and authoritarian regimes have sought to erode its foundations. Yet, despite these challenges, the spirit of democracy endures, testament to its resilience and the unwavering faith humanity has in the power of collective self-governance.'

$block2New = 'Continuing along the interwoven threads of governance, we encounter the concepts of authority, legitimacy, and accountability. Authority, the power to make decisions and enforce them, is a fundamental element of any governance system. However, authority is not simply bestowed upon those in power; it must be earned and maintained through legitimacy. Legitimacy, the perceived right to exercise authority, is essential for the stability and effectiveness of governance. When citizens perceive the government as legitimate, they are more likely to comply with its laws and policies. Accountability, closely intertwined with legitimacy, ensures that those in power are answerable for their actions. Through mechanisms such as elections, public scrutiny, and the rule of law, citizens hold their leaders accountable for their decisions. In the delicate balancing act of governance, accountability serves as a check against the potential abuse of power.'

$block3New = 'Finally, as we delve deeper into the tapestry of governance, we arrive at the intersection of ethics and decision-making. Ethics, a set of moral principles that guide human behavior, play a pivotal role in shaping the character of governance. Ethical decision-making requires leaders to weigh the potential consequences of their actions, considering the impact on both individuals and society as a whole. The choices made by those in power have far-reaching implications, affecting the well-being of citizens, the environment, and future generations. It is through ethical decision-making that leaders strive to uphold the common good, promote justice, and ensure the sustainability of our societies. In the intricate dance of governance, ethics serves as a compass, guiding leaders toward decisions that align with the values and aspirations of the people they serve.'

$summaryNew = 'In the labyrinthine realm of governance, where power and influence are intricately intertwined, diverse threads weave together to form a tapestry of democratic principles, authority, legitimacy, accountability, and ethics. Democracy, the foundation upon which modern governance rests, empowers citizens to shape policies and hold their leaders accountable. Authority, coupled with legitimacy and accountability, ensures the stability and effectiveness of governance. Ethical decision-making serves as a moral compass, guiding leaders toward choices that uphold the common good and promote justice. As we navigate the complexities of governance, it is imperative to recognize the interconnectedness of these elements and their profound impact on the well-being of our societies.'

$bodyNew = $block1New + $vtab + $vtab + $block2New + $vtab + $vtab + $block3New

# --- Replace the whole essay-body paragraph (paragraph 5) ---
$bodyPara = $d.Paragraphs(5)
$bodyStart = $bodyPara.Range.Start
$bodyEnd = $bodyPara.Range.End
$bodyRange = $d.Range($bodyStart, $bodyEnd)
$bodyRange.Text = $bodyNew

# --- Replace the Summary paragraph (paragraph 7) ---
$sumPara = $d.Paragraphs(7)
$sumStart = $sumPara.Range.Start
$sumEnd = $sumPara.Range.End
$sumRange = $d.Range($sumStart, $sumEnd)
$sumRange.Text = $summaryNew

# --- Append an extra empty paragraph at the end of the document ---
$d.Paragraphs.Add() | Out-Null

Write-Host "done"
